$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.5322491666365465
$ws.Cells.Item(2, 3).Value2 = 0.1489415781938419
$ws.Cells.Item(2, 4).Value2 = 0.04285421344128082
$ws.Cells.Item(2, 5).Value2 = 0.4085644444610779
$ws.Cells.Item(2, 6).Value2 = 1.326576608400103
$ws.Cells.Item(2, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(2, 9).Value2 = 0.9006924609901432
$ws.Cells.Item(2, 11).Value2 = 0.6678272695218368
$ws.Cells.Item(2, 14).Value2 = 1.833030973039698

$ws.Cells.Item(3, 2).Value2 = 0.4838363878315874
$ws.Cells.Item(3, 3).Value2 = 0.1334357997532152
$ws.Cells.Item(3, 4).Value2 = 0.04310138727917856
$ws.Cells.Item(3, 5).Value2 = 0.3563942191804301
$ws.Cells.Item(3, 6).Value2 = 1.303922196792854
$ws.Cells.Item(3, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(3, 9).Value2 = 0.8946564304235238
$ws.Cells.Item(3, 11).Value2 = 0.6046018203161907
$ws.Cells.Item(3, 14).Value2 = 1.844413391541984

$ws.Cells.Item(4, 2).Value2 = 0.4543624287547061
$ws.Cells.Item(4, 3).Value2 = 0.1239830996396734
$ws.Cells.Item(4, 4).Value2 = 0.04325566482178989
$ws.Cells.Item(4, 5).Value2 = 0.3244724349155348
$ws.Cells.Item(4, 6).Value2 = 1.290810902224848
$ws.Cells.Item(4, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(4, 9).Value2 = 0.8914314546030013
$ws.Cells.Item(4, 11).Value2 = 0.5660951004766162
$ws.Cells.Item(4, 14).Value2 = 1.852011909364464

$ws.Cells.Item(5, 2).Value2 = 0.4424144269603971
$ws.Cells.Item(5, 3).Value2 = 0.1201477732939793
$ws.Cells.Item(5, 4).Value2 = 0.04331916815619952
$ws.Cells.Item(5, 5).Value2 = 0.3114897092474109
$ws.Cells.Item(5, 6).Value2 = 1.285667909521365
$ws.Cells.Item(5, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(5, 9).Value2 = 0.8902378216162035
$ws.Cells.Item(5, 11).Value2 = 0.5504814614170073
$ws.Cells.Item(5, 14).Value2 = 1.855261267521534

$ws.Cells.Item(6, 2).Value2 = 0.4404342561145143
$ws.Cells.Item(6, 3).Value2 = 0.1195119202905346
$ws.Cells.Item(6, 4).Value2 = 0.04332975120342297
$ws.Cells.Item(6, 5).Value2 = 0.3093354223704381
$ws.Cells.Item(6, 6).Value2 = 1.284825972328449
$ws.Cells.Item(6, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(6, 9).Value2 = 0.8900468896575759
$ws.Cells.Item(6, 11).Value2 = 0.5478935266655469
$ws.Cells.Item(6, 14).Value2 = 1.855810044003213

$ws.Cells.Item(7, 2).Value2 = 0.4542010397815375
$ws.Cells.Item(7, 3).Value2 = 0.1239313078258704
$ws.Cells.Item(7, 4).Value2 = 0.04325651868078673
$ws.Cells.Item(7, 5).Value2 = 0.3242972443614462
$ws.Cells.Item(7, 6).Value2 = 1.290740733503497
$ws.Cells.Item(7, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(7, 9).Value2 = 0.8914148692330528
$ws.Cells.Item(7, 11).Value2 = 0.5658842140353499
$ws.Cells.Item(7, 14).Value2 = 1.852055112811797

$ws.Cells.Item(8, 2).Value2 = 0.5155039376268178
$ws.Cells.Item(8, 3).Value2 = 0.1435808740338302
$ws.Cells.Item(8, 4).Value2 = 0.04293891901645441
$ws.Cells.Item(8, 5).Value2 = 0.3905517577016724
$ws.Cells.Item(8, 6).Value2 = 1.318599092568803
$ws.Cells.Item(8, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(8, 9).Value2 = 0.8985111103451615
$ws.Cells.Item(8, 11).Value2 = 0.645961475150699
$ws.Cells.Item(8, 14).Value2 = 1.836828849818197

$ws.Cells.Item(9, 2).Value2 = 0.6377411981188743
$ws.Cells.Item(9, 3).Value2 = 0.1826703475786928
$ws.Cells.Item(9, 4).Value2 = 0.04233592847687717
$ws.Cells.Item(9, 5).Value2 = 0.5214642814763266
$ws.Cells.Item(9, 6).Value2 = 1.379608874740811
$ws.Cells.Item(9, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(9, 9).Value2 = 0.91626571309871
$ws.Cells.Item(9, 11).Value2 = 0.8055281737803739
$ws.Cells.Item(9, 14).Value2 = 1.811825996508432

$ws.Cells.Item(10, 2).Value2 = 0.7288270885583756
$ws.Cells.Item(10, 3).Value2 = 0.2117565391003211
$ws.Cells.Item(10, 4).Value2 = 0.0419048426121611
$ws.Cells.Item(10, 5).Value2 = 0.6184056745955644
$ws.Cells.Item(10, 6).Value2 = 1.428388380223979
$ws.Cells.Item(10, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(10, 9).Value2 = 0.9316821761194163
$ws.Cells.Item(10, 11).Value2 = 0.9243815531401083
$ws.Cells.Item(10, 14).Value2 = 1.796442202010823

$ws.Cells.Item(11, 2).Value2 = 0.7705524616292223
$ws.Cells.Item(11, 3).Value2 = 0.225074548946111
$ws.Cells.Item(11, 4).Value2 = 0.04171129564396914
$ws.Cells.Item(11, 5).Value2 = 0.6627068190713601
$ws.Cells.Item(11, 6).Value2 = 1.451453276089879
$ws.Cells.Item(11, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(11, 9).Value2 = 0.9392177479239479
$ws.Cells.Item(11, 11).Value2 = 0.9788196776068503
$ws.Cells.Item(11, 14).Value2 = 1.790097694592887

$ws.Cells.Item(12, 2).Value2 = 0.7863951258755719
$ws.Cells.Item(12, 3).Value2 = 0.230130628851299
$ws.Cells.Item(12, 4).Value2 = 0.04163837139436311
$ws.Cells.Item(12, 5).Value2 = 0.6795142584739011
$ws.Cells.Item(12, 6).Value2 = 1.460314252152386
$ws.Cells.Item(12, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(12, 9).Value2 = 0.9421469536610942
$ws.Cells.Item(12, 11).Value2 = 0.9994884738221401
$ws.Cells.Item(12, 14).Value2 = 1.787789673025358

$ws.Cells.Item(13, 2).Value2 = 0.7829812413673949
$ws.Cells.Item(13, 3).Value2 = 0.2290411342740128
$ws.Cells.Item(13, 4).Value2 = 0.04165406057078069
$ws.Cells.Item(13, 5).Value2 = 0.6758930318469396
$ws.Cells.Item(13, 6).Value2 = 1.458400226884237
$ws.Cells.Item(13, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(13, 9).Value2 = 0.94151272416903
$ws.Cells.Item(13, 11).Value2 = 0.9950346511191128
$ws.Cells.Item(13, 14).Value2 = 1.78828253538353

$ws.Cells.Item(14, 2).Value2 = 0.771854999203299
$ws.Cells.Item(14, 3).Value2 = 0.225490256111442
$ws.Cells.Item(14, 4).Value2 = 0.04170528874826029
$ws.Cells.Item(14, 5).Value2 = 0.6640889302757529
$ws.Cells.Item(14, 6).Value2 = 1.452179726887422
$ws.Cells.Item(14, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(14, 9).Value2 = 0.9394572164113413
$ws.Cells.Item(14, 11).Value2 = 0.9805190192857367
$ws.Cells.Item(14, 14).Value2 = 1.789905915284692

$ws.Cells.Item(15, 2).Value2 = 0.7650453609034571
$ws.Cells.Item(15, 3).Value2 = 0.2233169244081239
$ws.Cells.Item(15, 4).Value2 = 0.04173671539699875
$ws.Cells.Item(15, 5).Value2 = 0.6568627619730592
$ws.Cells.Item(15, 6).Value2 = 1.448386034050287
$ws.Cells.Item(15, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(15, 9).Value2 = 0.9382080253555216
$ws.Cells.Item(15, 11).Value2 = 0.9716348686935987
$ws.Cells.Item(15, 14).Value2 = 1.790912605057358

$ws.Cells.Item(16, 2).Value2 = 0.7261061122512444
$ws.Cells.Item(16, 3).Value2 = 0.2108879509363248
$ws.Cells.Item(16, 4).Value2 = 0.04191754280596705
$ws.Cells.Item(16, 5).Value2 = 0.6155147692807503
$ws.Cells.Item(16, 6).Value2 = 1.426898724295597
$ws.Cells.Item(16, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(16, 9).Value2 = 0.9312002668923895
$ws.Cells.Item(16, 11).Value2 = 0.9208314375738098
$ws.Cells.Item(16, 14).Value2 = 1.79687002722315

$ws.Cells.Item(17, 2).Value2 = 0.702292741991613
$ws.Cells.Item(17, 3).Value2 = 0.2032856243754395
$ws.Cells.Item(17, 4).Value2 = 0.04202912926982361
$ws.Cells.Item(17, 5).Value2 = 0.5902026297404603
$ws.Cells.Item(17, 6).Value2 = 1.413941732251487
$ws.Cells.Item(17, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(17, 9).Value2 = 0.9270354288371649
$ws.Cells.Item(17, 11).Value2 = 0.8897608618557058
$ws.Cells.Item(17, 14).Value2 = 1.800692491028158

$ws.Cells.Item(18, 2).Value2 = 0.6886231829302574
$ws.Cells.Item(18, 3).Value2 = 0.1989210904336574
$ws.Cells.Item(18, 4).Value2 = 0.04209355140873594
$ws.Cells.Item(18, 5).Value2 = 0.5756626194183525
$ws.Cells.Item(18, 6).Value2 = 1.406571498063172
$ws.Cells.Item(18, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(18, 9).Value2 = 0.9246890778696013
$ws.Cells.Item(18, 11).Value2 = 0.8719247661937857
$ws.Cells.Item(18, 14).Value2 = 1.802952565514758

$ws.Cells.Item(19, 2).Value2 = 0.6839995729816906
$ws.Cells.Item(19, 3).Value2 = 0.1974447175531679
$ws.Cells.Item(19, 4).Value2 = 0.04211540495472832
$ws.Cells.Item(19, 5).Value2 = 0.5707427856591494
$ws.Cells.Item(19, 6).Value2 = 1.404090165742616
$ws.Cells.Item(19, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(19, 9).Value2 = 0.9239030706131217
$ws.Cells.Item(19, 11).Value2 = 0.86589173404397
$ws.Cells.Item(19, 14).Value2 = 1.803728334181017

$ws.Cells.Item(20, 2).Value2 = 0.7048248910092809
$ws.Cells.Item(20, 3).Value2 = 0.2040940614680267
$ws.Cells.Item(20, 4).Value2 = 0.04201722580507727
$ws.Cells.Item(20, 5).Value2 = 0.5928951798732243
$ws.Cells.Item(20, 6).Value2 = 1.415312502843079
$ws.Cells.Item(20, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(20, 9).Value2 = 0.9274736921075615
$ws.Cells.Item(20, 11).Value2 = 0.8930647615576675
$ws.Cells.Item(20, 14).Value2 = 1.80027921516087

$ws.Cells.Item(21, 2).Value2 = 0.7751218978560814
$ws.Cells.Item(21, 3).Value2 = 0.2265328835134994
$ws.Cells.Item(21, 4).Value2 = 0.04169023180497788
$ws.Cells.Item(21, 5).Value2 = 0.6675552044231807
$ws.Cells.Item(21, 6).Value2 = 1.454003389171191
$ws.Cells.Item(21, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(21, 9).Value2 = 0.9400589117188929
$ws.Cells.Item(21, 11).Value2 = 0.9847811324193572
$ws.Cells.Item(21, 14).Value2 = 1.789426520291656

$ws.Cells.Item(22, 2).Value2 = 0.8213110909340458
$ws.Cells.Item(22, 3).Value2 = 0.2412729701817113
$ws.Cells.Item(22, 4).Value2 = 0.04147866647890108
$ws.Cells.Item(22, 5).Value2 = 0.7165349768227856
$ws.Cells.Item(22, 6).Value2 = 1.480029683914495
$ws.Cells.Item(22, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(22, 9).Value2 = 0.9487252564709223
$ws.Cells.Item(22, 11).Value2 = 1.045039778079058
$ws.Cells.Item(22, 14).Value2 = 1.782884750697917

$ws.Cells.Item(23, 2).Value2 = 0.7966363761136677
$ws.Cells.Item(23, 3).Value2 = 0.2333989132148133
$ws.Cells.Item(23, 4).Value2 = 0.04159138649898164
$ws.Cells.Item(23, 5).Value2 = 0.6903757731969336
$ws.Cells.Item(23, 6).Value2 = 1.466070957124074
$ws.Cells.Item(23, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(23, 9).Value2 = 0.9440593281664462
$ws.Cells.Item(23, 11).Value2 = 1.012849327666032
$ws.Cells.Item(23, 14).Value2 = 1.786325620449375

$ws.Cells.Item(24, 2).Value2 = 0.7036800408634747
$ws.Cells.Item(24, 3).Value2 = 0.2037285479508171
$ws.Cells.Item(24, 4).Value2 = 0.04202260652059486
$ws.Cells.Item(24, 5).Value2 = 0.5916778399362954
$ws.Cells.Item(24, 6).Value2 = 1.414692531722181
$ws.Cells.Item(24, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(24, 9).Value2 = 0.9272754035829394
$ws.Cells.Item(24, 11).Value2 = 0.891570985120552
$ws.Cells.Item(24, 14).Value2 = 1.800465862557914

$ws.Cells.Item(25, 2).Value2 = 0.6044516289907449
$ws.Cells.Item(25, 3).Value2 = 0.1720331571383724
$ws.Cells.Item(25, 4).Value2 = 0.04249694873450682
$ws.Cells.Item(25, 5).Value2 = 0.4859270393275779
$ws.Cells.Item(25, 6).Value2 = 1.362414310711102
$ws.Cells.Item(25, 8).Value2 = 0.07973214163530429
$ws.Cells.Item(25, 9).Value2 = 0.9110483733830392
$ws.Cells.Item(25, 11).Value2 = 0.6678272695218368
$ws.Cells.Item(25, 14).Value2 = 1.818067466884813
